$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Replace placeholder project data (title/description/skills) with the
#    real portfolio project content. The "image" column (E) keeps the same
#    value except for row 7 which becomes "blogpost-1".
# ---------------------------------------------------------------------------

$ws.Range("B2").Value = "Movie Recommender System with Sentiment Analysis"
$ws.Range("C2").Value = "Content Based Recommender System recommends movies similar to the movie user likes, and analyses the sentiments on the reviews given by the user for that movie."
$ws.Range("D2").Value = "HTML, CSS, JavaScript, Python, Bootstrap, Beautiful Soup, Scikit-Learn, Flask, Heroku"

$ws.Range("B3").Value = "Portfolio Website "
$ws.Range("C3").Value = "This project showcases my data science projects on a deployed website."
$ws.Range("D3").Value = "HTML, CSS, Python, Bootstrap, Flask, PostgreSQL, Heroku"

$ws.Range("B4").Value = "Human Rights First Asylum"
$ws.Range("C4").Value = "HRF needs a web tool backed by data science to aggregate data on asylum cases, allow users to explore that data, and predict and visualize how a judge might rule on a specific asylum case as well as what specific elements of an asylum case seem to most impact a favorable or unfavorable ruling"
$ws.Range("D4").Value = "Python, FastAPI,  AWS RDS PostgreSQL, AWS Elastic Beanstalk, Scikit-Learn, Docker, Selenium"

$ws.Range("B5").Value = "Airbnb in Los Angeles"
$ws.Range("C5").Value = "Predicting the nightly rates of Airbnb in Los Angeles depending on the location, property type, number of bedrooms, etc."
$ws.Range("D5").Value = "HTML, CSS, Python, PlotlyDash, Flask, plotly, Supervised Machine Learning"

$ws.Range("B6").Value = "Image Classifier using VGG-19 CNN"
$ws.Range("C6").Value = "Image Classification using Keras VGG-19 transfer learning"
$ws.Range("D6").Value = "HTML, CSS, JavaScript, Python, Flask, TensorFlow, Keras"

$ws.Range("B7").Value = "Airbus, the new King of the Skies?!"
$ws.Range("C7").Value = "Data Storytelling - Airbus v Boeing - How do they compare?"
$ws.Range("D7").Value = "Python, Data Wrangling, Data Visualization"
$ws.Range("E7").Value = "blogpost-1"

# ---------------------------------------------------------------------------
# 2. Move the "predicting airbnb prices" blog hyperlink from G3 to G5
#    (keeping the same target URL and the hyperlink/"Hyperlink" styling).
# ---------------------------------------------------------------------------

# Cut preserves both the cell value and the cell style (G3 ends up blank
# but keeps its Hyperlink style, matching the target workbook).
$ws.Range("G3").Cut($ws.Range("G5")) | Out-Null

# The old Hyperlink object still references G3 (now empty) - remove it and
# recreate it pointing at the new location, then restore the Hyperlink
# style (Hyperlinks.Add forces its own style xf onto the cell).
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$G$3') {
        $hl.Delete()
    }
}
$ws.Hyperlinks.Add($ws.Range("G5"), "https://navrozlamba.com/2020/08/28/predicting-airbnb-prices-in-los-angeles/") | Out-Null
$ws.Range("G5").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 3. Widen columns B, C, D and G to fit the new (longer) text.
#    ColumnWidth (chars) differs from the OOXML stored width by a constant
#    offset in this engine, so we subtract it to land on the exact target.
# ---------------------------------------------------------------------------

$offset = 0.8333333333333357
$ws.Columns.Item(2).ColumnWidth = 46.5 - $offset
$ws.Columns.Item(3).ColumnWidth = 59.5 - $offset
$ws.Columns.Item(4).ColumnWidth = 63 - $offset
$ws.Columns.Item(7).ColumnWidth = 67.5 - $offset

# ---------------------------------------------------------------------------
# 4. Update the active cell selection shown in the sheet view.
# ---------------------------------------------------------------------------

$ws.Range("E8").Select() | Out-Null
